# Auto-generated script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("D2", '54.332.46')
    ,@("E2", '  +0.38%  ')
    ,@("D3", '2.267.15')
    ,@("E3", '  -0.65%  ')
    ,@("E4", '  +0.16%  ')
    ,@("D5", '496.64')
    ,@("E5", '  +0.37%  ')
    ,@("D6", '128.91')
    ,@("E6", '  +1.00%  ')
    ,@("E7", '  -0.16%  ')
    ,@("E8", '  -0.51%  ')
    ,@("D9", '0.0955')
    ,@("E9", '  +0.75%  ')
    ,@("E10", '  +1.12%  ')
    ,@("D11", '0.338')
    ,@("E11", '  +3.69%  ')
    ,@("D12", '4.80')
    ,@("E12", '  +3.54%  ')
    ,@("D13", '22.96')
    ,@("E13", '  +5.23%  ')
    ,@("D14", '2.669.52')
    ,@("E14", '  -0.53%  ')
    ,@("D15", '54.297.44')
    ,@("E15", '  +0.33%  ')
    ,@("E16", '  +0.26%  ')
    ,@("D17", '2.277.24')
    ,@("E17", '  -0.52%  ')
    ,@("E18", '  +1.96%  ')
    ,@("E19", '  +1.37%  ')
    ,@("D20", '302.09')
    ,@("E20", '  +0.53%  ')
    ,@("D21", '6.34')
    ,@("E22", '  +0.08%  ')
    ,@("D23", '60.90')
    ,@("E23", '  -2.53%  ')
    ,@("D24", '0.997')
    ,@("E24", '  -1.07%  ')
    ,@("E25", '  +1.21%  ')
    ,@("E26", '  +3.86%  ')
    ,@("D27", '171.26')
    ,@("E27", '  +1.14%  ')
    ,@("B28", 'Aptos')
    ,@("C28", 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt')
    ,@("D28", '5.99')
    ,@("E28", '  +2.40%  ')
    ,@("D29", '1.61')
    ,@("E29", '  +0.26%  ')
    ,@("B30", 'PEPE')
    ,@("C30", 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe')
    ,@("D30", "0.0$([char]0x2083)0692")
    ,@("E30", '  +0.69%  ')
    ,@("E31", '  +1.79%  ')
    ,@("E32", '  +0.01%  ')
    ,@("D33", '17.83')
    ,@("E33", '  +1.04%  ')
    ,@("E34", '  -0.13%  ')
    ,@("D35", '0.937')
    ,@("E35", '  +8.77%  ')
    ,@("D36", '1.20')
    ,@("E37", '  +0.54%  ')
    ,@("E38", '  +0.23%  ')
    ,@("E39", '  +0.04%  ')
    ,@("D40", '3.38')
    ,@("E40", '  +0.91%  ')
    ,@("D41", '125.34')
    ,@("E41", '  -1.58%  ')
    ,@("D42", '4.81')
    ,@("E42", '  +0.80%  ')
    ,@("D43", '0.0494')
    ,@("E43", '  +2.32%  ')
    ,@("D44", '0.0894')
    ,@("E44", '  +0.82%  ')
    ,@("E45", '  +0.96%  ')
    ,@("D46", '241.67')
    ,@("E46", '  +1.93%  ')
    ,@("D47", '0.374')
    ,@("E47", '  +0.49%  ')
    ,@("E48", '  +1.50%  ')
    ,@("D49", '10.79')
    ,@("E49", '  +0.64%  ')
    ,@("D50", '16.11')
    ,@("E50", '  -0.65%  ')
    ,@("B51", 'BitgetToken')
    ,@("C51", 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb')
    ,@("D51", '0.934')
    ,@("E51", '  -0.53%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    if ($newVal -match '^[+-]?[0-9]+\.[0-9]+$' -or $newVal -match '^[+-]?[0-9]+$') {
        $ws.Range($cellRef).Value = "'" + $newVal
        $ws.Range($cellRef).Style = "Normal"
    } else {
        $ws.Range($cellRef).Value = $newVal
    }
}
